# Updates following review at IARC meeting 26/6
#
# - Move the "please see general notes" hints on the Submission sheet from
#   column C (C5/C6) to column B (B5/B6).
# - Make "Submission" the active tab (it was "Notes"), with B5:B6 selected.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Submission")

# Shift the two shared-string cells one column to the left.
$ws1.Range("B5").Value = $ws1.Range("C5").Value2
$ws1.Range("B6").Value = $ws1.Range("C6").Value2
$ws1.Range("C5").ClearContents()
$ws1.Range("C6").ClearContents()

# Switch the active tab from "Notes" back to "Submission", selecting B5:B6.
$ws1.Activate()
$ws1.Range("B5:B6").Select()
